$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Julio de 2020 a las 20:35"

# 2. Swap the country names for Fiyi / Dominica (rows 205 and 206) while keeping
#    their existing statistics (B:H) untouched, matching the shared-string
#    reorder in the diff (Fiyi now appears before Dominica).
$ws.Range("A205").Value = "Fiyi"
$ws.Range("A206").Value = "Dominica"

# 3. Update statistics for the updated countries.

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2964797
$ws.Range("C4").Value = 29027
$ws.Range("D4").Value = 1273950
$ws.Range("E4").Value = 1558374
$ws.Range("G4").Value = 155
$ws.Range("H4").Value = 132473

# Row 23 - Canada
$ws.Range("B23").Value = 105533
$ws.Range("C23").Value = 216
$ws.Range("D23").Value = 69239
$ws.Range("E23").Value = 27610
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = 8684

# Row 65 - Marruecos
$ws.Range("B65").Value = 14215
$ws.Range("C65").Value = 393
$ws.Range("D65").Value = 9725
$ws.Range("E65").Value = 4255
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 235

# Row 71 - Uzbekistan
$ws.Range("B71").Value = 9936
$ws.Range("C71").Value = 228
$ws.Range("D71").Value = 6446
$ws.Range("E71").Value = 3457
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 33

# Row 93 - Republica de Yibuti
$ws.Range("B93").Value = 4792
$ws.Range("C93").Value = 56
$ws.Range("D93").Value = 4593
$ws.Range("E93").Value = 144

# Row 108 - Maldivas
$ws.Range("B108").Value = 2468
$ws.Range("C108").Value = 33
$ws.Range("D108").Value = 2049
$ws.Range("E108").Value = 408
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 11

# Row 128 - Yemen
$ws.Range("B128").Value = 1265
$ws.Range("C128").Value = 17
$ws.Range("D128").Value = 552
$ws.Range("E128").Value = 375
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 338

# Row 130 - Tunez
$ws.Range("B130").Value = 1188
$ws.Range("C130").Value = 2
$ws.Range("E130").Value = 92

# Row 180 - Monaco
$ws.Range("B180").Value = 108
$ws.Range("C180").Value = 2
$ws.Range("E180").Value = 9
